$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Inputs!D13 ("Datamart") switches from "N" to "Y" so the Datamart
#    product is now priced under the old-pricing model.
# ------------------------------------------------------------------
$wsInputs = $wb.Worksheets.Item("Inputs")
$wsInputs.Range("D13").Value = "Y"

# ------------------------------------------------------------------
# 2. Outputs_Internal!D82 gains a guard so the Datamart line is
#    excluded from the discount-split totals whenever "Enable Old
#    Pricing" (Inputs!C4) is turned on ("Yes").
# ------------------------------------------------------------------
$wsInternal = $wb.Worksheets.Item("Outputs_Internal")
$wsInternal.Range("D82").Formula = '=IF(Inputs!$C$4 = "Yes", 0, SUM(Calcs!D74:H74)+SUM(Calcs!D85:H85))'

# ------------------------------------------------------------------
# 3. Leave the cursor where the author ended up after making the
#    edits, matching the saved selection/navigation state.
# ------------------------------------------------------------------
$wsInputs.Activate()
$wsInputs.Range("C4").Select()

$wsInternal.Activate()
$wsInternal.Range("D82").Select()

$wsTimeline = $wb.Worksheets.Item("Outputs_Timeline")
$wsTimeline.Activate()
$wsTimeline.Range("R27:R29").Select()
